# "maj export stock htp engagement"
#
# The stock-HTP export rows (Feuil1, rows 11-20) were filled in with a
# generic "/" placeholder for the B (heure) and I (heure) columns. This
# update replaces that placeholder with the real processing-time markers:
# "/17h/" for the stock-total rows (11-12) and "/23h/" for the "tnontj"
# rows (13-20) - mirroring the pattern already used for the other B/I
# columns further up the sheet.
#
# It also reflects the fact that the workbook was last saved with Feuil1
# selected (rather than Feuil2), updating each sheet's tabSelected flag
# and remembered selection accordingly.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Feuil1")
$ws2 = $wb.Worksheets.Item("Feuil2")

# --- Feuil1 : rows 11-12 -> "/17h/", rows 13-20 -> "/23h/" -----------------
$ws1.Range("B11").Value = "/17h/"
$ws1.Range("I11").Value = "/17h/"
$ws1.Range("B12").Value = "/17h/"
$ws1.Range("I12").Value = "/17h/"

$ws1.Range("B13").Value = "/23h/"
$ws1.Range("I13").Value = "/23h/"
$ws1.Range("B14").Value = "/23h/"
$ws1.Range("I14").Value = "/23h/"
$ws1.Range("B15").Value = "/23h/"
$ws1.Range("I15").Value = "/23h/"
$ws1.Range("B16").Value = "/23h/"
$ws1.Range("I16").Value = "/23h/"
$ws1.Range("B17").Value = "/23h/"
$ws1.Range("I17").Value = "/23h/"
$ws1.Range("B18").Value = "/23h/"
$ws1.Range("I18").Value = "/23h/"
$ws1.Range("B19").Value = "/23h/"
$ws1.Range("I19").Value = "/23h/"
$ws1.Range("B20").Value = "/23h/"
$ws1.Range("I20").Value = "/23h/"

# --- selection / active-sheet bookkeeping ----------------------------------
# Feuil2 keeps its own remembered selection (no longer the active tab).
$ws2.Range("J19").Select()

# Feuil1 becomes the active sheet/tab, with its own remembered selection.
$ws1.Activate()
$ws1.Range("G26").Select()
